# Add the "wl" worksheet to the workbook and append the new data/priors,
# then restore the previously active sheet ("Bm1") and its selection.

$wb = $excel.ActiveWorkbook

# --- Create the new "wl" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wl = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wl.Name = "wl"

# --- Fill in the id/val table for the new parameter (wl) ---
$rows = @(
  @("id", "val"),
  @("MNWA", 165),
  @("MWAP1", 91),
  @("MWAP1", 112),
  @("MWAP1", 640),
  @("MBATS", 13),
  @("MBATS", 90),
  @("MBATS", 23),
  @("MBATS", 70),
  @("MSARG", 9.74),
  @("MSARG", 12.9),
  @("MSARG", 29.3),
  @("MSARG", 15.7),
  @("MSARG", 19),
  @("MSARG", 36),
  @("MSARG", 45.5),
  @("MSARG", 54.1),
  @("MSARG", 8.83),
  @("MSARG", 11.2),
  @("MSARG", 19.600000000000001),
  @("MSARG", 7.7),
  @("MSARG", 8.7200000000000006),
  @("MSARG", 16.2),
  @("MSARG", 7.02),
  @("MSARG", 9.9600000000000009),
  @("MSARG", 15.9),
  @("MSARG", 5.81),
  @("MSARG", 7.75),
  @("MSARG", 14.6),
  @("MSARG", 17.600000000000001),
  @("MSARG", 9.58),
  @("MSARG", 5.23),
  @("MSARG", 7.65),
  @("MSARG", 11.1),
  @("MSARG", 10.199999999999999),
  @("MSARG", 10.1),
  @("MSARG", 17),
  @("MSARG", 21.9),
  @("MSARG", 13.6),
  @("MSARG", 5.12),
  @("MSARG", 2.5099999999999998),
  @("MSARG", 11.7),
  @("MSARG", 14.7),
  @("MSARG", 7.97),
  @("MSARG", 5.9),
  @("MSARG", 13.2),
  @("MSARG", 16.3),
  @("MSARG", 6.82),
  @("MSARG", 12.8),
  @("MSARG", 29.1),
  @("MSARG", 2.93),
  @("MSARG", 17.5),
  @("MSARG", 11.5),
  @("MSARG", 25.1),
  @("MWAP2", 78.787878800000001),
  @("MWAP2", 24.242424199999999),
  @("MWAP2", 22.222222200000001),
  @("MWAP2", 30.3030303),
  @("MWAP2", 47.474747499999999),
  @("MWAP2", 53.535353499999999),
  @("MWAP2", 143.43434300000001),
  @("MWAP2", 66.666666699999993),
  @("MWAP2", 38.383838400000002),
  @("MWAP2", 35.353535399999998),
  @("MWAP2", 73.737373700000006),
  @("MWAP2", 80.808080799999999),
  @("MWAP2", 91.919191900000001),
  @("MWAP2", 57.575757600000003),
  @("MWAP2", 38.383838400000002),
  @("MWAP2", 45.454545500000002),
  @("MWAP2", 72.7272727),
  @("MWAP2", 76.767676800000004),
  @("MWAP2", 149.49494899999999),
  @("MWAP2", 87.878787900000006),
  @("MWAP2", 44.444444400000002)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $wl.Cells.Item($r, 1).Value = $rows[$i][0]
    $wl.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Leave the cursor where it ended up on "wl" while editing (matches source edit)
$wl.Range("K14").Select()

# --- Restore focus to "Bm1", scrolled/selected near the bottom where rows were added ---
$bm1 = $wb.Worksheets.Item("Bm1")
$bm1.Activate()
$bm1.Range("A45:B45").Select()
